$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Replace the "X"/"O" text markers with real Boolean TRUE/FALSE values
# X -> TRUE, O -> FALSE
$ws.Range("B3").Value = $true
$ws.Range("C3").Value = $false
$ws.Range("D3").Value = $true

$ws.Range("B4").Value = $false
$ws.Range("C4").Value = $true
$ws.Range("D4").Value = $false

$ws.Range("B5").Value = $false
$ws.Range("C5").Value = $true
$ws.Range("D5").Value = $false

# Update the active selection to match the target state
$ws.Range("C6").Select()
